# Scheduled runner update: refresh market-price derived columns (H:N) across all
# job tables (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) with the latest Universalis pricing.
# Only numeric value cells change; row/column layout and styling are untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 80
$ws.Range("H80").Value = 2662.0527
$ws.Range("I80").Value = 3050.6667
$ws.Range("J80").Value = 2482.6924
$ws.Range("K80").Value = 9152.000100000001
$ws.Range("L80").Value = 7448.0772
$ws.Range("M80").Value = -8154.000100000001
$ws.Range("N80").Value = -9444.0772

# Row 83
$ws.Range("H83").Value = 2662.0527
$ws.Range("I83").Value = 3050.6667
$ws.Range("J83").Value = 2482.6924
$ws.Range("K83").Value = 27456.0003
$ws.Range("L83").Value = 22344.2316
$ws.Range("M83").Value = -22464.0003
$ws.Range("N83").Value = -32328.2316

# Row 86
$ws.Range("H86").Value = 3096.5
$ws.Range("I86").Value = 2923.5715
$ws.Range("J86").Value = 3500
$ws.Range("K86").Value = 2923.5715
$ws.Range("L86").Value = 3500
$ws.Range("M86").Value = -1800.5715
$ws.Range("N86").Value = -5746

# Row 89
$ws.Range("H89").Value = 3096.5
$ws.Range("I89").Value = 2923.5715
$ws.Range("J89").Value = 3500
$ws.Range("K89").Value = 14617.8575
$ws.Range("L89").Value = 17500
$ws.Range("M89").Value = -9001.8575
$ws.Range("N89").Value = -28732

# Row 107
$ws.Range("H107").Value = 4456.2856
$ws.Range("I107").Value = 4000
$ws.Range("J107").Value = 4532.3335
$ws.Range("K107").Value = 4000
$ws.Range("L107").Value = 4532.3335
$ws.Range("M107").Value = -2080
$ws.Range("N107").Value = -8372.333500000001

# Row 132
$ws.Range("H132").Value = 4775.549
$ws.Range("I132").Value = 5434.372
$ws.Range("J132").Value = 1234.375
$ws.Range("K132").Value = 16303.116
$ws.Range("L132").Value = 3703.125
$ws.Range("M132").Value = -13773.116
$ws.Range("N132").Value = -8763.125

# Row 137
$ws.Range("H137").Value = 10208.479
$ws.Range("I137").Value = 10482.527
$ws.Range("J137").Value = 9221.9
$ws.Range("K137").Value = 31447.581
$ws.Range("L137").Value = 27665.7
$ws.Range("M137").Value = -28897.581
$ws.Range("N137").Value = -32765.7

# Row 138
$ws.Range("H138").Value = 23811548
$ws.Range("I138").Value = 1224.24
$ws.Range("J138").Value = 58826730
$ws.Range("K138").Value = 3672.72
$ws.Range("L138").Value = 176480190
$ws.Range("M138").Value = 1467.28
$ws.Range("N138").Value = -176490470

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5411.341
$ws.Range("I32").Value = 5430.9287
$ws.Range("K32").Value = 5430.9287
$ws.Range("M32").Value = -5143.9287

# Row 61
$ws.Range("H61").Value = 3304.1
$ws.Range("I61").Value = 2894.16
$ws.Range("K61").Value = 2894.16
$ws.Range("M61").Value = -2682.16

# Row 74
$ws.Range("H74").Value = 48895.19
$ws.Range("I74").Value = 50811.64
$ws.Range("J74").Value = 984
$ws.Range("K74").Value = 50811.64
$ws.Range("L74").Value = 984
$ws.Range("M74").Value = -49937.64
$ws.Range("N74").Value = -2732

# Row 77
$ws.Range("H77").Value = 48895.19
$ws.Range("I77").Value = 50811.64
$ws.Range("J77").Value = 984
$ws.Range("K77").Value = 254058.2
$ws.Range("L77").Value = 4920
$ws.Range("M77").Value = -249690.2
$ws.Range("N77").Value = -13656

# Row 97
$ws.Range("H97").Value = 4986.1
$ws.Range("I97").Value = 4982.625
$ws.Range("K97").Value = 4982.625
$ws.Range("M97").Value = -4486.625

# Row 136
$ws.Range("H136").Value = 3304.1
$ws.Range("I136").Value = 2894.16
$ws.Range("K136").Value = 8682.48
$ws.Range("M136").Value = -6132.48

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 4254.3335
$ws.Range("I20").Value = 3315
$ws.Range("J20").Value = 6133
$ws.Range("K20").Value = 3315
$ws.Range("L20").Value = 6133
$ws.Range("M20").Value = -3068
$ws.Range("N20").Value = -6627

# Row 94
$ws.Range("H94").Value = 3047.5
$ws.Range("J94").Value = 2949.375
$ws.Range("L94").Value = 2949.375
$ws.Range("N94").Value = -3851.375

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1310
$ws.Range("I16").Value = 1200
$ws.Range("K16").Value = 1200
$ws.Range("M16").Value = -913

# Row 31
$ws.Range("H31").Value = 5614.879
$ws.Range("I31").Value = 3650.8823
$ws.Range("J31").Value = 7701.625
$ws.Range("K31").Value = 3650.8823
$ws.Range("L31").Value = 7701.625
$ws.Range("M31").Value = -3355.8823
$ws.Range("N31").Value = -8291.625

# Row 34
$ws.Range("H34").Value = 5614.879
$ws.Range("I34").Value = 3650.8823
$ws.Range("J34").Value = 7701.625
$ws.Range("K34").Value = 3650.8823
$ws.Range("L34").Value = 7701.625
$ws.Range("M34").Value = -3448.8823
$ws.Range("N34").Value = -8105.625

# Row 44
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()

# Row 86
$ws.Range("H86").Value = 10894.9
$ws.Range("I86").Value = 9668.666999999999
$ws.Range("J86").Value = 11420.429
$ws.Range("K86").Value = 9668.666999999999
$ws.Range("L86").Value = 11420.429
$ws.Range("M86").Value = -8545.666999999999
$ws.Range("N86").Value = -13666.429

# Row 89
$ws.Range("H89").Value = 10894.9
$ws.Range("I89").Value = 9668.666999999999
$ws.Range("J89").Value = 11420.429
$ws.Range("K89").Value = 48343.335
$ws.Range("L89").Value = 57102.145
$ws.Range("M89").Value = -42727.335
$ws.Range("N89").Value = -68334.145

# Row 105
$ws.Range("H105").Value = 18572.818
$ws.Range("I105").Value = 12448.571
$ws.Range("K105").Value = 12448.571
$ws.Range("M105").Value = -10701.571

# Row 113
$ws.Range("H113").Value = 1310
$ws.Range("I113").Value = 1200
$ws.Range("K113").Value = 1200
$ws.Range("M113").Value = 970

# Row 122
$ws.Range("H122").Value = 3668.1667
$ws.Range("I122").Value = 2901.8
$ws.Range("K122").Value = 8705.400000000001
$ws.Range("M122").Value = -6255.400000000001

$ws = $wb.Worksheets.Item("CUL")
# Row 17
$ws.Range("H17").Value = 350
$ws.Range("I17").Value = 366.66666
$ws.Range("K17").Value = 1099.99998
$ws.Range("M17").Value = -930.9999800000001

# Row 121
$ws.Range("H121").Value = 710.9091
$ws.Range("J121").Value = 957.1429000000001
$ws.Range("L121").Value = 2871.4287
$ws.Range("N121").Value = -5491.4287

# Row 128
$ws.Range("H128").Value = 160000
$ws.Range("I128").Value = 160000
$ws.Range("K128").Value = 480000
$ws.Range("M128").Value = -475020

# Row 132
$ws.Range("H132").Value = 3152.3845
$ws.Range("I132").Value = 1559.4
$ws.Range("K132").Value = 14034.6
$ws.Range("M132").Value = -11504.6

# Row 133
$ws.Range("H133").Value = 4564
$ws.Range("J133").Value = 7274.5
$ws.Range("L133").Value = 21823.5
$ws.Range("N133").Value = -31943.5

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5999.6
$ws.Range("I70").Value = 5999.6
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 5999.6
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -5729.6
$ws.Range("N70").ClearContents()

# Row 73
$ws.Range("H73").Value = 5999.6
$ws.Range("I73").Value = 5999.6
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 5999.6
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -5063.6
$ws.Range("N73").ClearContents()

# Row 109
$ws.Range("H109").Value = 37122.5
$ws.Range("I109").Value = 31250
$ws.Range("J109").Value = 42995
$ws.Range("K109").Value = 31250
$ws.Range("L109").Value = 42995
$ws.Range("M109").Value = -30210
$ws.Range("N109").Value = -45075

# Row 122
$ws.Range("H122").Value = 1723.3
$ws.Range("I122").Value = 1499.3529
$ws.Range("J122").Value = 2992.3333
$ws.Range("K122").Value = 4498.0587
$ws.Range("L122").Value = 8976.999899999999
$ws.Range("M122").Value = -2048.0587
$ws.Range("N122").Value = -13876.9999

# Row 126
$ws.Range("H126").Value = 2083
$ws.Range("I126").Value = 1980.7333
$ws.Range("J126").Value = 2850
$ws.Range("K126").Value = 5942.199900000001
$ws.Range("L126").Value = 8550
$ws.Range("M126").Value = -3472.199900000001
$ws.Range("N126").Value = -13490

# Row 136
$ws.Range("H136").Value = 59441.668
$ws.Range("J136").Value = 59441.668
$ws.Range("L136").Value = 178325.004
$ws.Range("N136").Value = -183425.004

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 2471.1333
$ws.Range("I16").Value = 2088
$ws.Range("J16").Value = 3237.4
$ws.Range("K16").Value = 2088
$ws.Range("L16").Value = 3237.4
$ws.Range("M16").Value = -1918
$ws.Range("N16").Value = -3577.4

# Row 22
$ws.Range("H22").Value = 1737.8889
$ws.Range("I22").Value = 1457.7273
$ws.Range("J22").Value = 2178.1428
$ws.Range("K22").Value = 1457.7273
$ws.Range("L22").Value = 2178.1428
$ws.Range("M22").Value = -1162.7273
$ws.Range("N22").Value = -2768.1428

# Row 27
$ws.Range("H27").Value = 1737.8889
$ws.Range("I27").Value = 1457.7273
$ws.Range("J27").Value = 2178.1428
$ws.Range("K27").Value = 1457.7273
$ws.Range("L27").Value = 2178.1428
$ws.Range("M27").Value = -1350.7273
$ws.Range("N27").Value = -2392.1428

# Row 55
$ws.Range("H55").Value = 1478.5883
$ws.Range("I55").Value = 1320.75
$ws.Range("J55").Value = 1527.1538
$ws.Range("K55").Value = 1320.75
$ws.Range("L55").Value = 1527.1538
$ws.Range("M55").Value = -1147.75
$ws.Range("N55").Value = -1873.1538

# Row 93
$ws.Range("H93").Value = 3366.5715
$ws.Range("I93").Value = 3366.5715
$ws.Range("K93").Value = 3366.5715
$ws.Range("M93").Value = -2118.5715

$ws = $wb.Worksheets.Item("WVR")
# Row 14
$ws.Range("H14").Value = 104.82353
$ws.Range("I14").Value = 106.07143
$ws.Range("J14").Value = 99
$ws.Range("K14").Value = 106.07143
$ws.Range("L14").Value = 99
$ws.Range("M14").Value = 61.92856999999999
$ws.Range("N14").Value = -435

# Row 136
$ws.Range("H136").Value = 6156.75
$ws.Range("I136").Value = 4532.75
$ws.Range("J136").Value = 9946.083000000001
$ws.Range("K136").Value = 13598.25
$ws.Range("L136").Value = 29838.249
$ws.Range("M136").Value = -11048.25
$ws.Range("N136").Value = -34938.249

